$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2023-09-09 -> 2023-09-10, serial 45178 -> 45179) for every data row (2..530).
$ws.Range("C2:C530").Value = 45179
